# AutoEval_Gétain_Maël.xlsx — fill in the student's self-evaluation
# (grades per criterion, remark, date, apprentice & teacher names)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the selector glyph next to the "chosen" column for every
# competency row (8-15): the middle "o" becomes "x".
$marker = "o`nx`no"
$ws.Range("E8:E15").Value = $marker

# Grades chosen for each competency (column G = "LARGEMENT ACQUIS",
# except row 15 where J = "SUFFISANT" was picked).
$ws.Range("G8").Value = 5.5
$ws.Range("G9").Value = 5.5
$ws.Range("G10").Value = 6
$ws.Range("G11").Value = 6
$ws.Range("G12").Value = 5.5
$ws.Range("G13").Value = 5.5
$ws.Range("G14").Value = 6
$ws.Range("J15").Value = 5

# Remark box
$ws.Range("B19").Value = "RAS"

# Date / Apprenti / Enseignant fields
$ws.Range("M19").Value = [DateTime]"2022-11-28"
$ws.Range("M21").Value = "Maël Gétain"
$ws.Range("M23").Value = "XCL"
